$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill A1:A23 with the updated producer-code list.
# Values are written in this specific order so that the shared-string
# table ends up built in the same order as a natural top-to-bottom
# data-entry pass interleaved with the original row staying last.
$ws.Cells.Item(23, 1).Value = "AG6304A44"
$ws.Cells.Item(22, 1).Value = "AG6304A136"
$ws.Cells.Item(21, 1).Value = "AG6304A129"
$ws.Cells.Item(20, 1).Value = "AG6304A85"
$ws.Cells.Item(19, 1).Value = "AG6304A140"
$ws.Cells.Item(18, 1).Value = "AG6304A106"
$ws.Cells.Item(17, 1).Value = "AG8258A1"
$ws.Cells.Item(16, 1).Value = "AG6304A36"
$ws.Cells.Item(10, 1).Value = "AG6304A109"
$ws.Cells.Item(12, 1).Value = "AG6304A38"
$ws.Cells.Item(11, 1).Value = "AG6304A125"
$ws.Cells.Item(14, 1).Value = "AG6304A70"
$ws.Cells.Item(15, 1).Value = "AG6304A55"
$ws.Cells.Item(13, 1).Value = "AG6304A49"
$ws.Cells.Item(9, 1).Value = "AG6304A143"
$ws.Cells.Item(8, 1).Value = "AG6304A133"
$ws.Cells.Item(7, 1).Value = "AG6304A29"
$ws.Cells.Item(6, 1).Value = "AG6304A112"
$ws.Cells.Item(5, 1).Value = "AG6304A132"
$ws.Cells.Item(4, 1).Value = "AG6304A35"
$ws.Cells.Item(3, 1).Value = "AG6304A116"
$ws.Cells.Item(1, 1).Value = "AG6304A51"
$ws.Cells.Item(2, 1).Value = "AG6304A91"

# Make sure every new row carries the same cell style (bordered,
# wrap-text, text-format) that row 1 already had.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Select the full used rows, matching the saved view state.
$ws.Range("1:24").Select() | Out-Null
